$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep a Text format so numeric-looking strings
# (e.g. "220.23", "1.699.04") are stored as text, matching the
# source inlineStr cells, instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.707.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -6.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.699.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.53%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5146"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -12.70%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.24"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06308"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07360"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.700.65"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.529"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5813"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.929.52"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008504"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.68"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -13.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.713.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.016"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -8.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "187.46"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -11.03%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -8.10%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.63"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.521"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1162"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.74"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.356"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05659"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.342"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.518"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.500"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.648"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.024"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6019"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.356"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.688"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01615"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.94%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.100.88"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8610"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -10.53%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.97"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.856.38"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000110"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.182"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05244"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4319"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.61%  "
